$wb = $excel.ActiveWorkbook

$sides  = $wb.Worksheets.Item("Sides")
$troops = $wb.Worksheets.Item("Sheet2")

# Rename Sheet2 -> Troops
$troops.Name = "Troops"

# Header row
$troops.Range("A1").Value = "Name"
$troops.Range("B1").Value = "Side"
$troops.Range("C1").Value = "Troop Count"
$troops.Range("D1").Value = "Damage Dealt Per Troop"
$troops.Range("E1").Value = "Hit Points Per Troop"

# Data rows
$troops.Range("A2").Value = "First Test Side Troops"
$troops.Range("B2").Value = "First Test Side"
$troops.Range("C2").Value = 1000
$troops.Range("D2").Value = 10
$troops.Range("E2").Value = 10

$troops.Range("A3").Value = "Second Test Side Troops"
$troops.Range("B3").Value = "Second Test Side"
$troops.Range("C3").Value = 10000
$troops.Range("D3").Value = 1
$troops.Range("E3").Value = 1

# Give the header row the same (bold+underline) look used for the "Name" header
# on the Sides sheet, by copying its format rather than re-deriving a style
# (keeps the shared style table untouched, same as the original workbook).
$sides.Range("A1").Copy()
$troops.Range("A1:E1").PasteSpecial(-4122)

# Column widths (as close as the engine's column-width rounding allows to the
# authored widths of 19.5546875 / 14.44140625 / 11.88671875 / 21 / 17.88671875)
$troops.Columns.Item(1).ColumnWidth = 18.666666666666668
$troops.Columns.Item(2).ColumnWidth = 13.666666666666666
$troops.Columns.Item(3).ColumnWidth = 11
$troops.Columns.Item(4).ColumnWidth = 20.166666666666668
$troops.Columns.Item(5).ColumnWidth = 17

# Page setup for the new sheet
$troops.PageSetup.Orientation = 1

# Sides sheet's own selection moves to A2 (do this while Sides is still the
# active sheet, before Troops is activated below)
[void]$sides.Range("A2").Select()

# Troops becomes the active sheet/tab, with F3 selected
$troops.Activate()
[void]$troops.Range("F3").Select()
